# Apply numeric "want-to-go" count updates (and one cancellation) across
# the four worksheets of the 广州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 577
$ws1.Range("F6").Value  = 370
$ws1.Range("F8").Value  = 152
$ws1.Range("F10").Value = 221
$ws1.Range("F11").Value = 5995
$ws1.Range("F12").Value = 57
$ws1.Range("F13").Value = 48
$ws1.Range("F14").Value = 494
$ws1.Range("F17").Value = 359
$ws1.Range("F21").Value = 710
$ws1.Range("F22").Value = 144
$ws1.Range("F24").Value = 311
$ws1.Range("F26").Value = 64
$ws1.Range("F27").Value = 1817
$ws1.Range("F28").Value = 480

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 267
# Event cancelled: rename and mark the lowest price as "not for sale"
$ws2.Range("C7").Value = "广州·动漫钢琴鬼才Kyle Xian互动演奏会（取消）"
$ws2.Range("G7").Value = "不可售"

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 243

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 243
$ws4.Range("F3").Value  = 577
$ws4.Range("F6").Value  = 741
$ws4.Range("F8").Value  = 370
$ws4.Range("F10").Value = 152
$ws4.Range("F12").Value = 221
$ws4.Range("F13").Value = 5995
$ws4.Range("F14").Value = 57
$ws4.Range("F15").Value = 48
$ws4.Range("F16").Value = 267
$ws4.Range("F17").Value = 494
$ws4.Range("F20").Value = 359
$ws4.Range("C27").Value = "广州·动漫钢琴鬼才Kyle Xian互动演奏会（取消）"
$ws4.Range("G27").Value = "不可售"
$ws4.Range("F28").Value = 710
$ws4.Range("F32").Value = 144
$ws4.Range("F34").Value = 311
$ws4.Range("F36").Value = 64
$ws4.Range("F37").Value = 1817
$ws4.Range("F38").Value = 480
